$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'325.11"
$ws.Range("E2").Value = "'-1.41%"
$ws.Range("F2").Value = "'7-2-2023"
$ws.Range("G2").Value = "'1"

$ws.Range("D3").Value = "'43.81"
$ws.Range("E3").Value = "'-2.89%"
$ws.Range("F3").Value = "'7-2-2023"
$ws.Range("G3").Value = "'1"

$ws.Range("D4").Value = "'5.538"
$ws.Range("E4").Value = "'-1.43%"
$ws.Range("F4").Value = "'7-2-2023"
$ws.Range("G4").Value = "'1"

$ws.Range("D5").Value = "'0.08009"
$ws.Range("E5").Value = "'-1.39%"
$ws.Range("F5").Value = "'7-2-2023"
$ws.Range("G5").Value = "'1"

$ws.Range("D6").Value = "'1.955"
$ws.Range("E6").Value = "'1.25%"
$ws.Range("F6").Value = "'7-2-2023"
$ws.Range("G6").Value = "'1"

$ws.Range("B7").Value = 'GateToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D7").Value = "'4.293"
$ws.Range("E7").Value = "'-0.91%"
$ws.Range("F7").Value = "'7-2-2023"
$ws.Range("G7").Value = "'1"

$ws.Range("B8").Value = 'BTSEToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D8").Value = "'2.569"
$ws.Range("E8").Value = "'-7.25%"
$ws.Range("F8").Value = "'7-2-2023"
$ws.Range("G8").Value = "'1"

$ws.Range("B9").Value = 'MXToken'
$ws.Range("C9").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D9").Value = "'0.9461"
$ws.Range("E9").Value = "'-1.16%"
$ws.Range("F9").Value = "'7-2-2023"
$ws.Range("G9").Value = "'1"

$ws.Range("B10").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C10").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D10").Value = "'0.1119"
$ws.Range("E10").Value = "'-5.06%"
$ws.Range("F10").Value = "'7-2-2023"
$ws.Range("G10").Value = "'1"

$ws.Range("B11").Value = 'WazirX'
$ws.Range("C11").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D11").Value = "'0.1834"
$ws.Range("E11").Value = "'-4.04%"
$ws.Range("F11").Value = "'7-2-2023"
$ws.Range("G11").Value = "'1"

$ws.Range("B12").Value = 'MCDex'
$ws.Range("C12").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D12").Value = "'11.81"
$ws.Range("E12").Value = "'35.94%"
$ws.Range("F12").Value = "'7-2-2023"
$ws.Range("G12").Value = "'1"

$ws.Range("B13").Value = 'MandalaExchangeToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D13").Value = "'0.09631"
$ws.Range("E13").Value = "'-2.32%"
$ws.Range("F13").Value = "'7-2-2023"
$ws.Range("G13").Value = "'1"

$ws.Range("B14").Value = 'BitrueCoin'
$ws.Range("C14").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D14").Value = "'0.04687"
$ws.Range("E14").Value = "'12.06%"
$ws.Range("F14").Value = "'7-2-2023"
$ws.Range("G14").Value = "'1"

$ws.Range("B15").Value = 'BitMartToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D15").Value = "'0.1066"
$ws.Range("E15").Value = "'-0.16%"
$ws.Range("F15").Value = "'7-2-2023"
$ws.Range("G15").Value = "'1"

$ws.Range("B16").Value = 'BitForexToken'
$ws.Range("C16").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D16").Value = "'0.001265"
$ws.Range("E16").Value = "'-1.19%"
$ws.Range("F16").Value = "'7-2-2023"
$ws.Range("G16").Value = "'1"

$ws.Range("B17").Value = 'CoinExToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D17").Value = "'0.04052"
$ws.Range("E17").Value = "'-7.22%"
$ws.Range("F17").Value = "'7-2-2023"
$ws.Range("G17").Value = "'1"

$ws.Range("B18").Value = 'TigerCash'
$ws.Range("C18").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D18").Value = "'0.005783"
$ws.Range("E18").Value = "'-4.56%"
$ws.Range("F18").Value = "'7-2-2023"
$ws.Range("G18").Value = "'1"

$ws.Range("B19").Value = 'HotbitToken'
$ws.Range("C19").Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$ws.Range("D19").Value = "'0.004304"
$ws.Range("E19").Value = "'-6.58%"
$ws.Range("F19").Value = "'7-2-2023"
$ws.Range("G19").Value = "'1"

$ws.Range("B20").Value = 'LEO'
$ws.Range("C20").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D20").Value = "'3.372"
$ws.Range("E20").Value = "'-5.85%"
$ws.Range("F20").Value = "'7-2-2023"
$ws.Range("G20").Value = "'1"

$ws.Range("D21").Value = "'0.3478"
$ws.Range("E21").Value = "'-0.19%"
$ws.Range("F21").Value = "'7-2-2023"
$ws.Range("G21").Value = "'1"

$ws.Range("D22").Value = "'0.1406"
$ws.Range("E22").Value = "'3.04%"
$ws.Range("F22").Value = "'7-2-2023"
$ws.Range("G22").Value = "'1"

$ws.Range("D23").Value = "'0.2545"
$ws.Range("E23").Value = "'-1.72%"
$ws.Range("F23").Value = "'7-2-2023"
$ws.Range("G23").Value = "'1"

$ws.Range("D24").Value = "'0.001241"
$ws.Range("E24").Value = "'0.06%"
$ws.Range("F24").Value = "'7-2-2023"
$ws.Range("G24").Value = "'1"

$ws.Range("D25").Value = "'0.0001190"
$ws.Range("E25").Value = "'-3.46%"
$ws.Range("F25").Value = "'7-2-2023"
$ws.Range("G25").Value = "'1"

$ws.Range("D26").Value = "'0.0003745"
$ws.Range("E26").Value = "'-6.26%"
$ws.Range("F26").Value = "'7-2-2023"
$ws.Range("G26").Value = "'1"

$ws.Range("F27").Value = "'7-2-2023"
$ws.Range("G27").Value = "'1"

$ws.Range("F28").Value = "'7-2-2023"
$ws.Range("G28").Value = "'1"

$ws.Range("F29").Value = "'7-2-2023"
$ws.Range("G29").Value = "'1"

$ws.Range("F30").Value = "'7-2-2023"
$ws.Range("G30").Value = "'1"

$ws.Range("F31").Value = "'7-2-2023"
$ws.Range("G31").Value = "'1"

$ws.Range("F32").Value = "'7-2-2023"
$ws.Range("G32").Value = "'1"

$ws.Range("F33").Value = "'7-2-2023"
$ws.Range("G33").Value = "'1"

$ws.Range("F34").Value = "'7-2-2023"
$ws.Range("G34").Value = "'1"

$ws.Range("F35").Value = "'7-2-2023"
$ws.Range("G35").Value = "'1"

$ws.Range("F36").Value = "'7-2-2023"
$ws.Range("G36").Value = "'1"

$ws.Range("F37").Value = "'7-2-2023"
$ws.Range("G37").Value = "'1"

$ws.Range("D38").Value = "'0.02513"
$ws.Range("E38").Value = "'-8.19%"
$ws.Range("F38").Value = "'7-2-2023"
$ws.Range("G38").Value = "'1"

$ws.Range("D39").Value = "'0.05506"
$ws.Range("E39").Value = "'-2.57%"
$ws.Range("F39").Value = "'7-2-2023"
$ws.Range("G39").Value = "'1"

$ws.Range("D40").Value = "'0.007520"
$ws.Range("E40").Value = "'-2.19%"
$ws.Range("F40").Value = "'7-2-2023"
$ws.Range("G40").Value = "'1"

$ws.Range("D41").Value = "'0.1388"
$ws.Range("E41").Value = "'-1.03%"
$ws.Range("F41").Value = "'7-2-2023"
$ws.Range("G41").Value = "'1"

$ws.Range("D42").Value = "'0.007350"
$ws.Range("E42").Value = "'-35.16%"
$ws.Range("F42").Value = "'7-2-2023"
$ws.Range("G42").Value = "'1"

$ws.Range("D43").Value = "'0.002016"
$ws.Range("E43").Value = "'-3.51%"
$ws.Range("F43").Value = "'7-2-2023"
$ws.Range("G43").Value = "'1"

$ws.Range("D44").Value = "'0.008294"
$ws.Range("E44").Value = "'-12.63%"
$ws.Range("F44").Value = "'7-2-2023"
$ws.Range("G44").Value = "'1"

$ws.Range("D45").Value = "'0.00007102"
$ws.Range("E45").Value = "'-0.62%"
$ws.Range("F45").Value = "'7-2-2023"
$ws.Range("G45").Value = "'1"

$ws.Range("D46").Value = "'0.00000000750"
$ws.Range("E46").Value = "'-0.13%"
$ws.Range("F46").Value = "'7-2-2023"
$ws.Range("G46").Value = "'1"

$ws.Range("E47").Value = "'1.25%"
$ws.Range("F47").Value = "'7-2-2023"
$ws.Range("G47").Value = "'1"

$ws.Range("D48").Value = "'0.004319"
$ws.Range("E48").Value = "'25.07%"
$ws.Range("F48").Value = "'7-2-2023"
$ws.Range("G48").Value = "'1"

$ws.Range("D49").Value = "'0.00002100"
$ws.Range("E49").Value = "'-0.13%"
$ws.Range("F49").Value = "'7-2-2023"
$ws.Range("G49").Value = "'1"

$ws.Range("D50").Value = "'0.0002000"
$ws.Range("E50").Value = "'-0.13%"
$ws.Range("F50").Value = "'7-2-2023"
$ws.Range("G50").Value = "'1"

$ws.Range("F51").Value = "'7-2-2023"
$ws.Range("G51").Value = "'1"

